# Apply parameter updates to the potential_preg_untrt sheet and
# re-point the active sheet/selection to match the authored session.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("SimParameters")
$wsUntrt = $wb.Worksheets.Item("potential_preg_untrt")

# Update the underlying parameter values (column C) on potential_preg_untrt.
# All other changed cells across this workbook (columns C/E on rows 50-58,
# 91-99 of this sheet, and the mirrored tables on potential_preg_trt,
# potential_preec_untrt/trt, postpreec_preg, etc.) are formulas that
# recompute automatically from these inputs.
$wsUntrt.Range("C9").Value = 0.05
$wsUntrt.Range("C10").Value = 0.02
$wsUntrt.Range("C11").Value = 0.02
$wsUntrt.Range("C13").Value = 0.005
$wsUntrt.Range("C14").Value = 0.004
$wsUntrt.Range("C15").Value = 0.004
$wsUntrt.Range("C16").Value = 0.004
$wsUntrt.Range("C17").Value = 0.004

# Force a full recalculation so every dependent formula (within this sheet
# and the other sheets that reference it) picks up the new parameter values.
$excel.CalculateFullRebuild()

# Move the active tab/selection from SimParameters to potential_preg_untrt,
# selecting C2:C21 there, matching the final authoring session state.
$wsUntrt.Activate()
$wsUntrt.Range("C2:C21").Select()
